# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2162
$ws1.Range("F4").Value = 44
$ws1.Range("F5").Value = 11384
$ws1.Range("F9").Value = 11333
$ws1.Range("F13").Value = 1742
$ws1.Range("F14").Value = 5655
$ws1.Range("F16").Value = 3477

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2162
$ws4.Range("F5").Value = 44
$ws4.Range("F7").Value = 11384
$ws4.Range("F11").Value = 11333
$ws4.Range("F15").Value = 1742
$ws4.Range("F17").Value = 5655
$ws4.Range("F19").Value = 3477
